$wb = $excel.ActiveWorkbook

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value = 1454.96
$ws.Cells.Item(129, 10).Value = 1824.1875
$ws.Cells.Item(129, 12).Value = 5472.5625
$ws.Cells.Item(129, 14).Value = -15472.5625

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 3217.956
$ws.Cells.Item(132, 9).Value = 2296.2922
$ws.Cells.Item(132, 10).Value = 5522.115
$ws.Cells.Item(132, 11).Value = 6888.8766
$ws.Cells.Item(132, 12).Value = 16566.345
$ws.Cells.Item(132, 13).Value = -4358.8766
$ws.Cells.Item(132, 14).Value = -21626.345

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2260.4
$ws.Cells.Item(137, 9).Value = 4151
$ws.Cells.Item(137, 11).Value = 12453
$ws.Cells.Item(137, 13).Value = -9903

# Sheet ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 220.2
$ws.Cells.Item(4, 9).Value = 167
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 167
$ws.Cells.Item(4, 12).Value = 300
$ws.Cells.Item(4, 13).Value = -51
$ws.Cells.Item(4, 14).Value = -532

# Sheet ARM row 6
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 2000
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 13).ClearContents()

# Sheet ARM row 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 2521
$ws.Cells.Item(26, 9).Value = 868.3333
$ws.Cells.Item(26, 10).Value = 5000
$ws.Cells.Item(26, 11).Value = 868.3333
$ws.Cells.Item(26, 12).Value = 5000
$ws.Cells.Item(26, 13).Value = -538.3333
$ws.Cells.Item(26, 14).Value = -5660

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6941.421
$ws.Cells.Item(32, 9).Value = 6456.8335
$ws.Cells.Item(32, 11).Value = 6456.8335
$ws.Cells.Item(32, 13).Value = -6169.8335

# Sheet ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 572.8570999999999
$ws.Cells.Item(97, 9).Value = 601.6667
$ws.Cells.Item(97, 11).Value = 601.6667
$ws.Cells.Item(97, 13).Value = -105.6667

# Sheet BSM row 15
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 10000
$ws.Cells.Item(15, 10).Value = 10000
$ws.Cells.Item(15, 12).Value = 10000
$ws.Cells.Item(15, 14).Value = -10454

# Sheet BSM row 33
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(33, 8).Value = 30000
$ws.Cells.Item(33, 9).Value = 30000
$ws.Cells.Item(33, 11).Value = 30000
$ws.Cells.Item(33, 13).Value = -29664

# Sheet BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2119.4211
$ws.Cells.Item(99, 9).Value = 1300.75
$ws.Cells.Item(99, 10).Value = 3522.8572
$ws.Cells.Item(99, 11).Value = 1300.75
$ws.Cells.Item(99, 12).Value = 3522.8572
$ws.Cells.Item(99, 13).Value = 197.25
$ws.Cells.Item(99, 14).Value = -6518.8572

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 5517.136
$ws.Cells.Item(134, 9).Value = 2614.2273
$ws.Cells.Item(134, 10).Value = 8420.046
$ws.Cells.Item(134, 11).Value = 7842.6819
$ws.Cells.Item(134, 12).Value = 25260.138
$ws.Cells.Item(134, 13).Value = -5307.6819
$ws.Cells.Item(134, 14).Value = -30330.138

# Sheet CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 36666830
$ws.Cells.Item(6, 9).Value = 36666830
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 36666830
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).ClearContents()
$ws.Cells.Item(6, 14).Value = -36666717

# Sheet CRP row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 21285
$ws.Cells.Item(12, 9).Value = 475
$ws.Cells.Item(12, 10).Value = 52500
$ws.Cells.Item(12, 11).Value = 475
$ws.Cells.Item(12, 12).Value = 52500
$ws.Cells.Item(12, 13).Value = -305
$ws.Cells.Item(12, 14).Value = -52840

# Sheet CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 14666.667
$ws.Cells.Item(41, 9).Value = 9500
$ws.Cells.Item(41, 11).Value = 9500
$ws.Cells.Item(41, 13).Value = -9072

# Sheet CRP row 47
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(47, 8).Value = 50000
$ws.Cells.Item(47, 10).Value = 50000
$ws.Cells.Item(47, 12).Value = 50000
$ws.Cells.Item(47, 14).Value = -51132

# Sheet CRP row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(69, 8).Value = 3121.3333
$ws.Cells.Item(69, 9).Value = 3121.3333
$ws.Cells.Item(69, 11).Value = 3121.3333
$ws.Cells.Item(69, 13).Value = -2372.3333

# Sheet CRP row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(72, 8).Value = 3121.3333
$ws.Cells.Item(72, 9).Value = 3121.3333
$ws.Cells.Item(72, 11).Value = 9363.999899999999
$ws.Cells.Item(72, 13).Value = -5619.999899999999

# Sheet CRP row 93
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(93, 8).Value = 6365.5454
$ws.Cells.Item(93, 9).Value = 3540.125
$ws.Cells.Item(93, 10).Value = 13900
$ws.Cells.Item(93, 11).Value = 3540.125
$ws.Cells.Item(93, 12).Value = 13900
$ws.Cells.Item(93, 13).Value = -1668.125
$ws.Cells.Item(93, 14).Value = -17644

# Sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2066.8948
$ws.Cells.Item(99, 9).Value = 2002.3846
$ws.Cells.Item(99, 10).Value = 2206.6667
$ws.Cells.Item(99, 11).Value = 2002.3846
$ws.Cells.Item(99, 12).Value = 2206.6667
$ws.Cells.Item(99, 13).Value = -504.3846000000001
$ws.Cells.Item(99, 14).Value = -5202.6667

# Sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2066.8948
$ws.Cells.Item(126, 9).Value = 2002.3846
$ws.Cells.Item(126, 10).Value = 2206.6667
$ws.Cells.Item(126, 11).Value = 6007.1538
$ws.Cells.Item(126, 12).Value = 6620.000100000001
$ws.Cells.Item(126, 13).Value = -3537.1538
$ws.Cells.Item(126, 14).Value = -11560.0001

# Sheet CUL row 102
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(102, 8).Value = 3350
$ws.Cells.Item(102, 10).Value = 3350
$ws.Cells.Item(102, 12).Value = 10050
$ws.Cells.Item(102, 14).Value = -14918

# Sheet GSM row 36
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36, 8).Value = 3176.8333
$ws.Cells.Item(36, 9).Value = 1214.2858
$ws.Cells.Item(36, 10).Value = 4425.727
$ws.Cells.Item(36, 11).Value = 1214.2858
$ws.Cells.Item(36, 12).Value = 4425.727
$ws.Cells.Item(36, 13).Value = -729.2858000000001
$ws.Cells.Item(36, 14).Value = -5395.727

# Sheet GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 6800
$ws.Cells.Item(43, 9).Value = 1000
$ws.Cells.Item(43, 10).Value = 10666.667
$ws.Cells.Item(43, 11).Value = 1000
$ws.Cells.Item(43, 12).Value = 10666.667
$ws.Cells.Item(43, 13).Value = -849
$ws.Cells.Item(43, 14).Value = -10968.667

# Sheet LTW row 19
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 1661
$ws.Cells.Item(19, 9).Value = 1501.5
$ws.Cells.Item(19, 10).Value = 1980
$ws.Cells.Item(19, 11).Value = 1501.5
$ws.Cells.Item(19, 12).Value = 1980
$ws.Cells.Item(19, 13).Value = -1331.5
$ws.Cells.Item(19, 14).Value = -2320

# Sheet LTW row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 30000
$ws.Cells.Item(20, 9).Value = 27777.777
$ws.Cells.Item(20, 10).Value = 36666.668
$ws.Cells.Item(20, 11).Value = 27777.777
$ws.Cells.Item(20, 12).Value = 36666.668
$ws.Cells.Item(20, 13).Value = -27551.777
$ws.Cells.Item(20, 14).Value = -37118.668

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5960
$ws.Cells.Item(122, 9).Value = 2800
$ws.Cells.Item(122, 11).Value = 8400
$ws.Cells.Item(122, 13).Value = -5950

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 27059114
$ws.Cells.Item(132, 9).Value = 33371774
$ws.Cells.Item(132, 11).Value = 100115322
$ws.Cells.Item(132, 13).Value = -100112792

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4305.5625
$ws.Cells.Item(132, 10).Value = 3545.2727
$ws.Cells.Item(132, 12).Value = 10635.8181
$ws.Cells.Item(132, 14).Value = -15695.8181

# Sheet WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 38361.5
$ws.Cells.Item(135, 10).Value = 38361.5
$ws.Cells.Item(135, 12).Value = 38361.5
$ws.Cells.Item(135, 14).Value = -48501.5

